# Update crypto price/volume figures per latest fetch.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the Price/Volume columns so values
# like "1.000" / "0.2860" are stored verbatim instead of being coerced
# to numbers (which would drop the meaningful trailing zeros).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.350.44"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "1.904.53"
$ws.Range("E3").Value = "  -2.79%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "238.73"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4779"
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("D8").Value = "0.2860"
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").Value = "0.06693"
$ws.Range("E9").Value = "  -4.49%  "
$ws.Range("D10").Value = "18.79"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "102.73"
$ws.Range("E11").Value = "  -4.45%  "
$ws.Range("D12").Value = "0.07718"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "1.913.73"
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").Value = "5.216"
$ws.Range("E14").Value = "  -4.86%  "
$ws.Range("D15").Value = "0.6753"
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").Value = "261.68"
$ws.Range("E16").Value = "  -7.52%  "
$ws.Range("D17").Value = "30.358.10"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "0.000007500"
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").Value = "12.73"
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("D21").Value = "5.438"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "6.298"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").Value = "9.475"
$ws.Range("E24").Value = "  -3.73%  "
$ws.Range("D25").Value = "164.31"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Value = "18.98"
$ws.Range("E26").Value = "  -5.05%  "
$ws.Range("D27").Value = "2.070"
$ws.Range("E27").Value = "  -5.60%  "
$ws.Range("D28").Value = "0.1012"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").Value = "1.378"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "4.630"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "1.514"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("D32").Value = "4.227"
$ws.Range("E32").Value = "  -4.39%  "
$ws.Range("D33").Value = "0.04773"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("D34").Value = "0.7328"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").Value = "1.115"
$ws.Range("E35").Value = "  -4.63%  "
$ws.Range("D36").Value = "0.9997"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "2.700"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").Value = "0.01922"
$ws.Range("E38").Value = "  -4.27%  "
$ws.Range("D39").Value = "2.587"
$ws.Range("E39").Value = "  -4.38%  "
$ws.Range("D40").Value = "6.292"
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").Value = "74.68"
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("D42").Value = "1.995"
$ws.Range("E42").Value = "  -5.64%  "
$ws.Range("D43").Value = "0.8615"
$ws.Range("E43").Value = "  -4.78%  "
$ws.Range("D44").Value = "106.33"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("D45").Value = "0.4257"
$ws.Range("E45").Value = "  -4.51%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "1.007.59"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("D48").Value = "7.477"
$ws.Range("E48").Value = "  -7.83%  "
$ws.Range("D49").Value = "35.01"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("D50").Value = "0.1193"
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("D51").Value = "8.827"
$ws.Range("E51").Value = "  -5.65%  "

# Restore the original (default/Normal) cell style now that the text
# values are safely in place, matching the workbook's original formatting.
$dataRange.Style = "Normal"
